$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.361.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.066.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.10"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.32%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0762"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.371.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.776"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.056.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.563.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +15.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "227.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.90%  "
$ws.Range("E28").Value = "  +13.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0619"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.80%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("E40").Value = "  +3.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0975"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.46%  "
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +22.21%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.56%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.456.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0211"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.94%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.90%  "
$ws.Range("E49").Value = "  +3.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.58%  "
$ws.Range("E51").Value = "  +2.05%  "
